$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "1 Message Sent" sheet: just move the selection (no longer active)
# ------------------------------------------------------------------
$wsMsg = $wb.Worksheets.Item("1 Message Sent")
$wsMsg.Range("H29:H31").Select() | Out-Null

# ------------------------------------------------------------------
# 2) "NAMED PIPE RESULTS" sheet: move the selection, drop tabSelected
#    (will be overridden later once another sheet is selected last)
# ------------------------------------------------------------------
$wsPipe = $wb.Worksheets.Item("NAMED PIPE RESULTS")
$wsPipe.Range("A3:A13").Select() | Out-Null

# ------------------------------------------------------------------
# 3) Add the new "WINSOCK RESULTS" sheet after "NAMED PIPE RESULTS"
#    (becomes sheetId 4 / rId4, last tab)
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsWinsock = $wb.Worksheets.Add($null, $lastSheet)
$wsWinsock.Name = "WINSOCK RESULTS"

# Copy the layout/formatting of the NAMED PIPE RESULTS sheet (same
# styles: s=7 title, s=4 header, s=2 run numbers, s=5 averages row)
$wsPipe.Range("A1:E13").Copy() | Out-Null
$wsWinsock.Range("A1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
# The source sheet has no A2 cell at all - PasteSpecial(formats) leaves
# a blank styleless cell there, so get rid of it again.
$wsWinsock.Range("A2").Clear() | Out-Null

# ------------------------------------------------------------------
# 4) Fill in the Winsock content/values
# ------------------------------------------------------------------
$wsWinsock.Range("A1").Value = "Winsock"

$wsWinsock.Range("B2").Value = "Run "
$wsWinsock.Range("C2").Value = "Time (ns)"
$wsWinsock.Range("D2").Value = "Average"
$wsWinsock.Range("E2").Value = "Named Pipes (40 Bytes)"

$data = @(
  @(2345909, 2366676, 3408044, 40499313),
  @(2816752, 2796741, 2419160, 36669511),
  @(2543384, 2245850, 2563396, 45141669),
  @(2550180, 2804670, 4356904, 40935419),
  @(2815242, 2389331, 2336848, 34904321),
  @(2815242, 2938334, 3912114, 37184909),
  @(2720847, 2783526, 2057438, 47468698),
  @(2717071, 4036716, 2755207, 46709384),
  @(2689885, 2857531, 3334794, 36363671),
  @(2598511, 3857365, 3908717, 32130990)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $row = 3 + $i
  $wsWinsock.Cells.Item($row, 1).Value = $i + 1
  $wsWinsock.Cells.Item($row, 2).Value = $data[$i][0]
  $wsWinsock.Cells.Item($row, 3).Value = $data[$i][1]
  $wsWinsock.Cells.Item($row, 4).Value = $data[$i][2]
  $wsWinsock.Cells.Item($row, 5).Value = $data[$i][3]
}

$wsWinsock.Range("A13").Value = "Average"
$wsWinsock.Range("B13").Formula = "=AVERAGE(B3:B12)"
$wsWinsock.Range("C13").Formula = "=AVERAGE(C3:C12)"
$wsWinsock.Range("D13").Formula = "=AVERAGE(D3:D12)"
$wsWinsock.Range("E13").Formula = "=AVERAGE(E3:E12)"

# ------------------------------------------------------------------
# 5) Final selection on the new sheet - this makes it the active tab
# ------------------------------------------------------------------
$wsWinsock.Range("I16").Select() | Out-Null
